# Update the "Denominadores_D1_119" export template:
#  - drop the unused Hoja2 / Hoja3 worksheets
#  - rename/re-order the Hoja1 header row (EN naming + new Age_Group /
#    ILI* / ICU* / Deaths* / Pneu* columns) and give every renamed header
#    cell the same bold/fill/date-number-format style already used by
#    column B (StartDateOfWeek)

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$ws = $wb.Worksheets.Item("Hoja1")

# --- remove the two empty placeholder sheets -------------------------------
$wb.Worksheets.Item("Hoja2").Delete()
$wb.Worksheets.Item("Hoja3").Delete()

# --- rewrite the header row (columns E..AD) ---------------------------------
$ws.Range("E1").Value  = "Age_Group"
$ws.Range("F1").Value  = "ILINumFem"
$ws.Range("G1").Value  = "ILINumMale"
$ws.Range("H1").Value  = "ILINumST"
$ws.Range("I1").Value  = "ILINumEmerST"
$ws.Range("J1").Value  = "ILIDenoFem"
$ws.Range("K1").Value  = "ILIDenoMale"
$ws.Range("L1").Value  = "ILIDenoST"
$ws.Range("M1").Value  = "HospFem"
$ws.Range("N1").Value  = "HospMale"
$ws.Range("O1").Value  = "HospST"
$ws.Range("P1").Value  = "ICUFem"
$ws.Range("Q1").Value  = "ICUMale"
$ws.Range("R1").Value  = "ICUST"
$ws.Range("S1").Value  = "DeathsFem"
$ws.Range("T1").Value  = "DeathsMale"
$ws.Range("U1").Value  = "DeathsST"
$ws.Range("V1").Value  = "PneuFem"
$ws.Range("W1").Value  = "PneuMale"
$ws.Range("X1").Value  = "PneuST"
$ws.Range("Y1").Value  = "CCSARIFem"
$ws.Range("Z1").Value  = "CCSARIMale"
$ws.Range("AA1").Value = "CCSARIST"
$ws.Range("AB1").Value = "VentFem"
$ws.Range("AC1").Value = "VentMale"
$ws.Range("AD1").Value = "VentST"

# --- match the header style used by column B (bold font, header fill, and
#     the yyyy/mm/dd custom number format) across the whole renamed range --
$ws.Range("E1:AD1").NumberFormat = "yyyy/mm/dd"
